$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "model" sheet: drop the trailing "sector" property row (old row 18).
#    It is no longer part of the model, so select then delete the whole row.
# ------------------------------------------------------------------
$modelWs = $wb.Worksheets.Item("model")
$modelWs.Rows.Item(18).Select() | Out-Null
$modelWs.Rows.Item(18).Delete() | Out-Null

# ------------------------------------------------------------------
# 2. "properties" sheet: add the two new "Table"/"security" rows that
#    configure the security behaviour for table creation.
# ------------------------------------------------------------------
$propsWs = $wb.Worksheets.Item("properties")

$propsWs.Range("A3").Value = "Table"
$propsWs.Range("B3").Value = "security"
$propsWs.Range("C3").Value = "unverifiedUserCanCreate"
$propsWs.Range("D3").Value = "boolean"
# leading apostrophe forces text storage so "false" stays a string value
# (matching the shared-string "false", not an Excel boolean TRUE/FALSE)
$propsWs.Range("E3").Value = "'false"

$propsWs.Range("A4").Value = "Table"
$propsWs.Range("B4").Value = "security"
$propsWs.Range("C4").Value = "filterTypeOnCreation"
$propsWs.Range("D4").Value = "string"
$propsWs.Range("E4").Value = "HIDDEN"

# ------------------------------------------------------------------
# 3. Make "properties" the active/visible tab with D7 selected, matching
#    where the author ended up after adding the new rows.
# ------------------------------------------------------------------
$propsWs.Activate() | Out-Null
$propsWs.Range("D7").Select() | Out-Null
